# Edit: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# - Add a new period (2507) to the account-statement detail table, keeping the
#   existing "last row" border styling on the bottom row of the table.
# - Re-order the period column to descending order (most recent period first).
# - Update the totals (VALOR MORA and Cant. Periodos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new detail row right before the current last data row (row 27),
#    so the existing bottom row (with its special bottom-border style) is
#    pushed down to row 28 and keeps that styling, while the newly inserted
#    row 27 takes on the regular interior-row styling (copied from row 26).
$ws.Rows("27:27").Insert()
$ws.Range("B26:J26").Copy()
$ws.Range("B27:J27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 2. Fill in the values for the newly inserted row (same worker/contribution
#    data as the other detail rows, new period).
$ws.Range("B27").Value2 = "CC"
$ws.Range("C27").Value2 = "1140868712"
$ws.Range("D27").Value2 = "HECTOR JULIO PEREZ MUÃ?OZ"
$ws.Range("E27").Value2 = "2412"
$ws.Range("F27").Value2 = 216110
$ws.Range("G27").Value2 = 5402745

# 3. Re-sort the period column (rows 16-28) in descending order and add the
#    new period 2507 at the top, pushing the oldest period down to row 28.
$periods = @("2507","2506","2505","2504","2503","2502","2501","2412","2411","2410","2409","2408","2407")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value2 = $periods[$i]
}

# 4. Update totals.
$ws.Range("E11").Value2 = 2809430
$ws.Range("F13").Value2 = 13

Write-Host "Edit applied"
